# ---------------------------------------------------------------------------
# Rename sheets and populate the (previously empty) second sheet with a
# "java_programs" table that mirrors the look of the "algo_programs" sheet,
# adding the first Java program (KargerMinCut).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the worksheets -----------------------------------------------
$wsAlgo = $wb.Worksheets.Item(1)
$wsJava = $wb.Worksheets.Item(2)
$wsAlgo.Name = "algo_programs"
$wsJava.Name = "java_programs"

# --- 2. Column widths for java_programs -------------------------------------
$wsJava.Columns.Item(2).ColumnWidth = 11.833333333333334   # -> ~12.71
$wsJava.Columns.Item(3).ColumnWidth = 30.166666666666668   # -> 31
$wsJava.Columns.Item(4).ColumnWidth = 14.0                 # -> ~14.86
$wsJava.Columns.Item(5).ColumnWidth = 93.66666666666667    # -> ~94.43
$wsJava.Columns.Item(6).ColumnWidth = 14.166666666666666   # -> 15

# --- 3. Cell values ----------------------------------------------------------
$wsJava.Cells.Item(2,2).Value = "Sr. No. "
$wsJava.Cells.Item(2,3).Value = "File Name"
$wsJava.Cells.Item(2,4).Value = "Language"
$wsJava.Cells.Item(2,5).Value = "Description"
$wsJava.Cells.Item(2,6).Value = "Status"

$wsJava.Cells.Item(3,2).Value = 1
$wsJava.Cells.Item(3,3).Value = "KargerMinCut"
$wsJava.Cells.Item(3,4).Value = "Java"
$wsJava.Cells.Item(3,5).Value = "Coursera - Stanford - Karger Min Cut Implementation - See description file for details"
$wsJava.Cells.Item(3,6).Value = "Completed"

# --- 4. Helper functions for styling -----------------------------------------
$White = 16777215

function Format-Border($cell, $hasRight, $bottomWeight) {
    if ($hasRight) {
        $rb = $cell.Borders.Item(10)
        $rb.LineStyle = 1
        $rb.Weight = 2
        $rb.Color = $White
    }
    $bb = $cell.Borders.Item(9)
    $bb.LineStyle = 1
    $bb.Weight = $bottomWeight
    $bb.Color = $White
}

function Format-HeaderCell($cell, $wrap) {
    $cell.Font.Bold = $true
    $cell.Font.Size = 12
    $cell.Font.Color = $White
    $cell.Interior.ThemeColor = 9   # theme accent5 (matches theme="8" in OOXML)
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4108    # xlCenter
    if ($wrap) { $cell.WrapText = $true }
}

function Format-BandCell($cell, $fillColor, $align, $wrap) {
    $cell.Font.Size = 12
    $cell.Interior.Color = $fillColor
    $cell.HorizontalAlignment = $align
    $cell.VerticalAlignment = -4108   # xlCenter
    if ($wrap) { $cell.WrapText = $true }
}

# Computed RGB equivalents of theme "accent5" (theme index 8) with the tints
# used by the original TableStyleMedium13-based formatting:
#   no tint            -> 4BACC6
#   tint 0.6 (lighter)  -> B7DEE8
#   tint 0.8 (lighter)  -> DBEEF4
$FillDark  = 13020235   # 0x4BACC6 (BGR-packed) - header
$FillMed   = 15261367   # 0xB7DEE8 (BGR-packed) - odd data rows
$FillLight = 16051931   # 0xDBEEF4 (BGR-packed) - even data rows

$xlCenter = -4108
$xlLeft = -4131

# --- 5. Header row (row 2) ---------------------------------------------------
foreach ($c in 2..4) {
    $cell = $wsJava.Cells.Item(2, $c)
    Format-HeaderCell $cell $false
    Format-Border $cell $true 4
}
$eCell = $wsJava.Cells.Item(2, 5)
Format-HeaderCell $eCell $true
Format-Border $eCell $true 4

$fCell = $wsJava.Cells.Item(2, 6)
Format-HeaderCell $fCell $false
Format-Border $fCell $false 4

# --- 6. Data / blank banded rows (rows 3-6) ----------------------------------
# Row 3 = data (darker band), Row 4 = blank (lighter band),
# Row 5 = blank (darker band), Row 6 = blank (lighter band)
$rowFills = @{ 3 = $FillMed; 4 = $FillLight; 5 = $FillMed; 6 = $FillLight }

foreach ($r in 3..6) {
    $fill = $rowFills[$r]
    foreach ($c in 2..4) {
        $cell = $wsJava.Cells.Item($r, $c)
        Format-BandCell $cell $fill $xlCenter $false
        Format-Border $cell $true 2
    }
    $eCell = $wsJava.Cells.Item($r, 5)
    Format-BandCell $eCell $fill $xlLeft $true
    Format-Border $eCell $true 2

    $fCell = $wsJava.Cells.Item($r, 6)
    Format-BandCell $fCell $fill $xlCenter $false
    Format-Border $fCell $false 2
}

# --- 7. Row heights -----------------------------------------------------------
$wsJava.Rows.Item(2).RowHeight = 16.5
$wsJava.Rows.Item(3).RowHeight = 16.5
$wsJava.Rows.Item(4).RowHeight = 15.75
$wsJava.Rows.Item(5).RowHeight = 15.75
$wsJava.Rows.Item(6).RowHeight = 15.75

# --- 8. View / selection state -------------------------------------------------
$wsAlgo.Activate()
$wsAlgo.Range("E16").Select()
$winAlgo = $excel.ActiveWindow
$winAlgo.ScrollRow = 19
$winAlgo.ScrollColumn = 1

$wsJava.Activate()
$wsJava.Range("E4").Select()

$wsAlgo.Activate()

Write-Output "Applied java_programs sheet edits"
